$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Mark the "CAMERA & VIEW" sub-rows (20, 31-36) as completed on Milestone II.
$rowsToMark = 20,31,32,33,34,35,36
foreach ($r in $rowsToMark) {
    $ws.Range("E$r").Value = "II"
    $ws.Range("F$r").Value = "X"
}

# Mark both "Effective Use of GIT" (row 90) and "All Graphics API Objects
# cleaned up in memory" (row 91) bonus rows as completed for Milestone I.
$ws.Range("C90").Value = "X"
$ws.Range("C91").Value = "X"

# Restore the on-screen selection to match the author's last edit location.
$ws.Range("C91").Select()
